$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.65
$ws.Range("H2").Value = 2.82
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1.52
$ws.Range("K2").Value = 2.22
$ws.Range("L2").Value = 2.47
$ws.Range("M2").Value = 1.42
$ws.Range("N2").Value = 1.6
$ws.Range("O2").Value = 2.07
$ws.Range("P2").Value = 2.05
$ws.Range("R2").Value = 6.3
$ws.Range("S2").Value = 11.75
$ws.Range("T2").Value = 10.5
$ws.Range("U2").Value = 32
$ws.Range("V2").Value = 28
$ws.Range("W2").Value = 50
$ws.Range("X2").Value = 6.1
$ws.Range("Y2").Value = 5.7
$ws.Range("AB2").Value = 6.5
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 35
$ws.Range("AF2").Value = 30
$ws.Range("AG2").Value = 50

# Row 5 updates
$ws.Range("G5").Value = 2.75
$ws.Range("H5").Value = 2.92
$ws.Range("I5").Value = 2.6
$ws.Range("K5").Value = 2.45
$ws.Range("L5").Value = 2.25
$ws.Range("M5").Value = 1.5
$ws.Range("O5").Value = 2.32
$ws.Range("P5").Value = 1.9
$ws.Range("Q5").Value = 1.72
$ws.Range("R5").Value = 7.3
$ws.Range("S5").Value = 13
$ws.Range("T5").Value = 10.25
$ws.Range("U5").Value = 32
$ws.Range("V5").Value = 26
$ws.Range("W5").Value = 40
$ws.Range("X5").Value = 7
$ws.Range("Y5").Value = 5.8
$ws.Range("Z5").Value = 16
$ws.Range("AA5").Value = 90
$ws.Range("AB5").Value = 6.7
$ws.Range("AC5").Value = 11.75
$ws.Range("AD5").Value = 10.25
$ws.Range("AE5").Value = 30
$ws.Range("AF5").Value = 26
$ws.Range("AG5").Value = 40

$wb.Save()
